{"js": "// Replace the single run-of-text Bibliografia paragraph (which concatenates\n// nine references back-to-back with no separators) with the same nine\n// references separated by manual line breaks (<w:br/>), matching how the\n// rest of the document already formats multi-line list-style paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the bibliography paragraph robustly (search by its distinctive\n// leading text) rather than relying on a fixed paragraph index.\nconst marker = \"DENNIS, P. (2009). Produ\u00e7\u00e3o lean simplificada.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Bibliografia paragraph not found\");\n}\n\n// The individual reference entries, in order, exactly as they appear in the\n// original run-on text (the last one keeps its trailing space before 1998).\nconst references = [\n  \"DENNIS, P. (2009). Produ\u00e7\u00e3o lean simplificada. Bookman Editora.\",\n  \"LEAN ENTERPRISE INSTITUTE (2007). L\u00e9xico lean \u2013 gloss\u00e1rio ilustrado para praticantes do pensamento lean. v.2.0. S\u00e3o Paulo: Lean Institute Brasil.\",\n  \"LIKER, J. K. (2005). O Modelo Toyota: 14 princ\u00edpios de gest\u00e3o do maior fabricante do mundo. Tradu\u00e7\u00e3o  de Lene Belon Ribeiro. Porto Alegre: Bookman.\",\n  \"ROSER, C. (2022). Tudo Sobre Produ\u00e7\u00e3o Puxada: Projetando, Implementando e Mantendo Kanban, CONWIP e outros Sistemas Puxados na Produ\u00e7\u00e3o Enxuta. AllAboutLean.com Publishing.\",\n  \"ROTHER, M.; HARRIS, R. (2002). Criando fluxo cont\u00ednuo. S\u00e3o Paulo, SP. Lean Institute Brasil.\",\n  \"SHOOK, John; ROTHER, Mike. Manual. Aprendendo a enxergar. Leam Institute Brasil. S\u00e3o Paulo: IMAM, s/d.\",\n  \"TUBINO, D. F. (2015). Manufatura enxuta como estrat\u00e9gia de produ\u00e7\u00e3o. Editora Atlas SA.\",\n  \"WOMACK, James P.; JONES, Daniel T. A Mentalidade enxuta nas empresas. Rio de Janeiro: Campus, 1998. \",\n  \"WOMACK, James P.; JONES, Daniel T. Lean Thinking: Banish Waste and Create Wealth in Your Corporation. Free Press, 2010.\"\n];\n\n// Office.js represents a Word manual line break (<w:br/>) as the vertical\n// tab character (\\v / \\u000b) inside paragraph/range text. Joining the\n// references with that character and writing it back via insertText\n// recreates the exact <w:t>/<w:br/> run structure from the diff.\nconst newText = references.join(\"\\u000b\");\n\ntarget.getRange().insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the single run-of-text Bibliografia paragraph (which concatenates\n# nine references back-to-back with no separators) with the same nine\n# references separated by manual line breaks, matching how the rest of the\n# document already formats multi-line list-style paragraphs.\n\n$d = $word.ActiveDocument\n\n# Locate the bibliography paragraph robustly (search by its distinctive\n# leading text) rather than relying on a fixed paragraph index.\n$marker = \"DENNIS, P. (2009). Produ\u00e7\u00e3o lean simplificada.\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith($marker)) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Bibliografia paragraph not found\"\n}\n\n# The individual reference entries, in order, exactly as they appear in the\n# original run-on text (the last one keeps its trailing space before 1998).\n$references = @(\n    \"DENNIS, P. (2009). Produ\u00e7\u00e3o lean simplificada. Bookman Editora.\",\n    \"LEAN ENTERPRISE INSTITUTE (2007). L\u00e9xico lean \u2013 gloss\u00e1rio ilustrado para praticantes do pensamento lean. v.2.0. S\u00e3o Paulo: Lean Institute Brasil.\",\n    \"LIKER, J. K. (2005). O Modelo Toyota: 14 princ\u00edpios de gest\u00e3o do maior fabricante do mundo. Tradu\u00e7\u00e3o  de Lene Belon Ribeiro. Porto Alegre: Bookman.\",\n    \"ROSER, C. (2022). Tudo Sobre Produ\u00e7\u00e3o Puxada: Projetando, Implementando e Mantendo Kanban, CONWIP e outros Sistemas Puxados na Produ\u00e7\u00e3o Enxuta. AllAboutLean.com Publishing.\",\n    \"ROTHER, M.; HARRIS, R. (2002). Criando fluxo cont\u00ednuo. S\u00e3o Paulo, SP. Lean Institute Brasil.\",\n    \"SHOOK, John; ROTHER, Mike. Manual. Aprendendo a enxergar. Leam Institute Brasil. S\u00e3o Paulo: IMAM, s/d.\",\n    \"TUBINO, D. F. (2015). Manufatura enxuta como estrat\u00e9gia de produ\u00e7\u00e3o. Editora Atlas SA.\",\n    \"WOMACK, James P.; JONES, Daniel T. A Mentalidade enxuta nas empresas. Rio de Janeiro: Campus, 1998. \",\n    \"WOMACK, James P.; JONES, Daniel T. Lean Thinking: Banish Waste and Create Wealth in Your Corporation. Free Press, 2010.\"\n)\n\n# Word represents a manual line break (<w:br/>) inside Range.Text as the\n# vertical-tab character Chr(11). Joining the references with that\n# character and writing it back to the range recreates the exact\n# <w:t>/<w:br/> run structure from the target document, without disturbing\n# the paragraph mark itself.\n$newText = [string]::Join([string][char]11, $references)\n\n$r = $target.Range\n$r.End = $r.End - 1   # exclude the trailing paragraph mark\n$r.Text = $newText\n"}
